# DS_Algo almost final version
# Adds the two new "Graph" rows to the Data sheet (new shared strings),
# matching the wrap-text style already used by the rest of the column,
# and reflows the affected row heights the same way Excel did when the
# sheet was resaved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New content rows (A34, A35) -> two new shared strings.
$ws.Cells.Item(34, 1).Value = 'print("I am from Graph link")'
$ws.Cells.Item(35, 1).Value = 'print("I am from Graph Representation")'

# Match the wrap-text style (s="1") already used by A3:A33.
$ws.Range("A34:A35").WrapText = $true

# Row-height reflow for the rows whose wrapped height changed.
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 120
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(29).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 30
$ws.Rows.Item(31).RowHeight = 30

# Move the selection/active cell down to the newly added last row.
$ws.Range("A35").Select()
